$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new data row (no special style)
$ws.Range("A2").Value = "images/adl-M-ver2.jpg"
$ws.Range("B2").Value = "M"
$ws.Range("C2").Value = "adolescent"

# Row 3: existing row (A3 already has style "1"), fill in values
$ws.Range("A3").Value = "images/adl-M-ver3.jpg"
$ws.Range("B3").Value = "M"
$ws.Range("C3").Value = "adolescent"

# Row 4: existing row (A4 already has style "1"), fill in values
$ws.Range("A4").Value = "images/adl-F-ver5.jpg"
$ws.Range("B4").Value = "F"
$ws.Range("C4").Value = "adolescent"

# Row 5: new data row (no special style)
$ws.Range("A5").Value = "images/adl-F-ver8.jpg"
$ws.Range("B5").Value = "F"
$ws.Range("C5").Value = "adolescent"

# Update selection to match final state
$ws.Range("C9").Select()

$wb.Save()
